# Trade #74 (MarketMaking) closes, and a new Trade #107 (MarketMaking) is opened.
# This touches the Summary, Strategy Status, All Trades and MarketMaking sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1400.49   # Current Capital
$wsSummary.Range("B4").Value = 0.29      # Total P&L $
$wsSummary.Range("B6").Value = 74        # Total Trades
$wsSummary.Range("B8").Value = 31        # Losing Trades
$wsSummary.Range("B9").Value = 44.59     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 100.49     # Capital
$wsStatus.Range("D5").Value = 41         # Trades
$wsStatus.Range("E5").Value = 0.18       # P&L $
$wsStatus.Range("F5").Value = 0.49       # P&L %
$wsStatus.Range("G5").Value = 46.34      # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet - update the closed trade (row 75, Trade #74)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")
$wsAll.Range("G75").Value = 0.119565
$wsAll.Range("H75").Value = "CLOSED"
$wsAll.Range("I75").Value = -8.0268
$wsAll.Range("J75").Value = -0.01
$wsAll.Range("K75").Value = 100.49
$wsAll.Range("L75").Value = "early_exit"
$wsAll.Range("M75").Value = 0.14

# Append the newly-opened trade (row 108, Trade #107)
$wsAll.Range("A108").Value = 107
# "2026-02-17" looks like a date to Excel's parser - force it to stay text by
# temporarily marking the cell as Text, then resetting the style so no
# stray number-format is left applied to the cell.
$wsAll.Range("B108").NumberFormat = "@"
$wsAll.Range("B108").Value = "2026-02-17"
$wsAll.Range("B108").Style = "Normal"
$wsAll.Range("C108").Value = "21:07:24"
$wsAll.Range("D108").Value = "MarketMaking"
$wsAll.Range("E108").Value = "DOWN"
$wsAll.Range("F108").Value = 0.13
# G108 stays blank (empty Exit Price, trade is still OPEN)
$wsAll.Range("H108").Value = "OPEN"
$wsAll.Range("I108").Value = 0
$wsAll.Range("J108").Value = 0
$wsAll.Range("K108").Value = 100.5019219857093
# L108 stays blank (empty Exit Reason, trade is still OPEN)
$wsAll.Range("M108").Value = 0
$wsAll.Range("N108").Value = 0
$wsAll.Range("O108").Value = 0
$wsAll.Range("P108").Value = 0.6
$wsAll.Range("Q108").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet - update the closed trade (row 42, Trade #74)
# ---------------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("G42").Value = 0.119565
$wsMM.Range("H42").Value = "CLOSED"
$wsMM.Range("I42").Value = -8.0268
$wsMM.Range("J42").Value = -0.01
$wsMM.Range("K42").Value = 100.49
$wsMM.Range("P42").Value = "early_exit"
$wsMM.Range("Q42").Value = 0.14

# Append the newly-opened trade (row 75, Trade #107)
$wsMM.Range("A75").Value = 107
$wsMM.Range("B75").NumberFormat = "@"
$wsMM.Range("B75").Value = "2026-02-17"
$wsMM.Range("B75").Style = "Normal"
$wsMM.Range("C75").Value = "21:07:24"
$wsMM.Range("D75").Value = "MarketMaking"
$wsMM.Range("E75").Value = "DOWN"
$wsMM.Range("F75").Value = 0.13
# G75 stays blank (empty Exit Price, trade is still OPEN)
$wsMM.Range("H75").Value = "OPEN"
$wsMM.Range("I75").Value = 0
$wsMM.Range("J75").Value = 0
$wsMM.Range("K75").Value = 100.5019219857093
$wsMM.Range("L75").Value = 0
$wsMM.Range("M75").Value = 0
$wsMM.Range("N75").Value = 0.6
$wsMM.Range("O75").Value = "Normal spread capture: 19600 bps"
# P75 stays blank (empty Exit Reason, trade is still OPEN)
$wsMM.Range("Q75").Value = 0
